$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-10 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-11 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("35÷5=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=5, 0", 2) | Out-Null
$d.Content.Find.Execute("16÷6=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "46÷4=11, 2", 2) | Out-Null
$d.Content.Find.Execute("20÷4=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "37÷6=6, 1", 2) | Out-Null
$d.Content.Find.Execute("44÷2=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷6=15, 3", 2) | Out-Null
$d.Content.Find.Execute("13÷2=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=5, 1", 2) | Out-Null
$d.Content.Find.Execute("61÷3=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "29÷2=14, 1", 2) | Out-Null
$d.Content.Find.Execute("66÷9=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2) | Out-Null
$d.Content.Find.Execute("79÷9=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=11, 7", 2) | Out-Null
$d.Content.Find.Execute("38÷2=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=22, 3", 2) | Out-Null
$d.Content.Find.Execute("17÷6=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=10, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=9, 3", 2) | Out-Null
$d.Content.Find.Execute("47÷8=5, 7", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("94÷4=23, 2", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2) | Out-Null
$d.Content.Find.Execute("60÷6=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "31÷7=4, 3", 2) | Out-Null
$d.Content.Find.Execute("61÷9=6, 7", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=11, 2", 2) | Out-Null
$d.Content.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 2) | Out-Null
$d.Content.Find.Execute("40÷7=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=6, 6", 2) | Out-Null
$d.Content.Find.Execute("64÷8=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=29, 1", 2) | Out-Null
$d.Content.Find.Execute("49÷5=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 2) | Out-Null
$d.Content.Find.Execute("31÷4=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=10, 1", 2) | Out-Null
$d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=10, 0", 2) | Out-Null
$d.Content.Find.Execute("76÷2=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=7, 0", 2) | Out-Null
$d.Content.Find.Execute("50÷4=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "61÷3=20, 1", 2) | Out-Null
$d.Content.Find.Execute("15÷7=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "15÷4=3, 3", 2) | Out-Null
$d.Content.Find.Execute("61÷6=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "87÷9=9, 6", 2) | Out-Null
